$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.402.83"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "1.730.32"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'322.33"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4543"
$ws.Range("E7").Value = "  +6.45%  "
$ws.Range("D8").Value = "'0.3521"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").Value = "'0.07349"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").Value = "'41.22"
$ws.Range("E10").Value = "  -3.13%  "
$ws.Range("D11").Value = "'1.074"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'20.38"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "'5.911"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("D15").Value = "'7.045"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").Value = "1.726.83"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").Value = "'91.05"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "'0.00001049"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").Value = "'0.06326"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'16.59"
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").Value = "'5.722"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("D23").Value = "27.447.76"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").Value = "'11.05"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "'2.063"
$ws.Range("E25").Value = "  -2.63%  "
$ws.Range("D26").Value = "'162.00"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'19.83"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").Value = "1.925.69"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("D29").Value = "'2.040"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("D30").Value = "'124.27"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").Value = "'1.040"
$ws.Range("E31").Value = "  -6.89%  "
$ws.Range("D32").Value = "'0.09109"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").Value = "'3.650"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").Value = "'5.349"
$ws.Range("E34").Value = "  -4.16%  "
$ws.Range("D35").Value = "'0.02260"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").Value = "'11.57"
$ws.Range("E36").Value = "  -5.50%  "
$ws.Range("D37").Value = "'0.05941"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").Value = "'0.2046"
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("D39").Value = "'0.6214"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").Value = "'4.848"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("D41").Value = "'1.189"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "'1.367"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").Value = "'7.688"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").Value = "'12.98"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "'3.694"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "'0.5781"
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("D47").Value = "'121.80"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("E48").Value = "  -3.55%  "
$ws.Range("D49").Value = "'0.06824"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'1.108"
$ws.Range("E50").Value = "  -6.06%  "
$ws.Range("D51").Value = "'70.85"
$ws.Range("E51").Value = "  -3.91%  "
